$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A15").Value = "особливый товар"
$ws.Range("A16").Value = "мелочь"
$ws.Range("A19").Value = "крамными товар"
$ws.Range("A20").Value = "небогатый товар"
$ws.Range("A24").Value = "щепетильный товар"
$ws.Range("A25").Value = "нужный товар"
$ws.Range("A26").Value = "набойчатый товар"
$ws.Range("A27").Value = "пушной товар"
$ws.Range("A28").Value = "медный товар"
$ws.Range("A29").Value = "недорогой товар"
$ws.Range("A31").Value = "питейный припасы"
$ws.Range("A32").Value = "внутренний товар"
$ws.Range("A35").Value = "заморский товар"
$ws.Range("A36").Value = "галантерейный товар"
$ws.Range("A37").Value = "произрастание"
$ws.Range("A38").Value = "купецкий товар"
$ws.Range("A39").Value = "домовый товар"
$ws.Range("A40").Value = "харчевой припасы"
$ws.Range("A41").Value = "рукодельный товар"
$ws.Range("A42").Value = "меховой товар"
$ws.Range("A43").Value = "надлежащий товар"
